# Regenerate save_data to use K (strikeouts thrown) instead of Strike#,
# regen std/mean, calc and write s_vals.
#
# The only observable change in this workbook's canonical OOXML is the
# value of column G ("K") for rows 2-11. Write the new K values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 3
    4  = 0
    5  = 1
    6  = 0
    7  = 2
    8  = 1
    9  = 1
    10 = 2
    11 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
